$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I33").Value = 521.3333
$ws.Range("H33").Value = 1007.38464
$ws.Range("K33").Value = 521.3333
$ws.Range("N33").Value = -1882
$ws.Range("L33").Value = 1424
$ws.Range("M33").Value = -292.3333
$ws.Range("J33").Value = 1424
$ws.Range("I62").Value = 7489.8237
$ws.Range("H62").Value = 7457.636
$ws.Range("J62").Value = 7348.2
$ws.Range("L62").Value = 7348.2
$ws.Range("K62").Value = 7489.8237
$ws.Range("M62").Value = -6865.8237
$ws.Range("N62").Value = -8596.200000000001
$ws.Range("H65").Value = 7457.636
$ws.Range("N65").Value = -42981
$ws.Range("L65").Value = 36741
$ws.Range("K65").Value = 37449.1185
$ws.Range("I65").Value = 7489.8237
$ws.Range("M65").Value = -34329.1185
$ws.Range("J65").Value = 7348.2
$ws.Range("M74").Value = -6172
$ws.Range("I74").Value = 7108
$ws.Range("H74").Value = 7614
$ws.Range("K74").Value = 7108
$ws.Range("M77").Value = -30860
$ws.Range("H77").Value = 7614
$ws.Range("I77").Value = 7108
$ws.Range("K77").Value = 35540
$ws.Range("M80").Value = -870.6999999999998
$ws.Range("H80").Value = 1018.76
$ws.Range("K80").Value = 1868.7
$ws.Range("I80").Value = 622.9
$ws.Range("M83").Value = -614.0999999999995
$ws.Range("H83").Value = 1018.76
$ws.Range("I83").Value = 622.9
$ws.Range("K83").Value = 5606.099999999999
$ws.Range("H98").Value = 1220.7778
$ws.Range("N98").Value = -6138.3333
$ws.Range("J98").Value = 3142.3333
$ws.Range("K98").Value = 980.5833
$ws.Range("L98").Value = 3142.3333
$ws.Range("M98").Value = 517.4167
$ws.Range("I98").Value = 980.5833
$ws.Range("I111").Value = 8800.35
$ws.Range("N111").Value = -19191.092
$ws.Range("M111").Value = -23334.05
$ws.Range("J111").Value = 4352.364
$ws.Range("K111").Value = 26401.05
$ws.Range("L111").Value = 13057.092
$ws.Range("H111").Value = 7222.032
$ws.Range("I116").Value = 2300
$ws.Range("M116").Value = 1142
$ws.Range("K116").Value = 2300
$ws.Range("H116").Value = 2300
$ws.Range("H122").Value = 1220.7778
$ws.Range("N122").Value = -14326.9999
$ws.Range("J122").Value = 3142.3333
$ws.Range("M122").Value = -491.7498999999998
$ws.Range("K122").Value = 2941.7499
$ws.Range("I122").Value = 980.5833
$ws.Range("L122").Value = 9426.999899999999
$ws.Range("M132").Value = -4259.3678
$ws.Range("H132").Value = 6683.407
$ws.Range("K132").Value = 6789.3678
$ws.Range("I132").Value = 2263.1226
$ws.Range("H137").Value = 2594.7
$ws.Range("L137").Value = 5621.3079
$ws.Range("N137").Value = -10721.3079
$ws.Range("J137").Value = 1873.7693
$ws.Range("H141").Value = 3049.4583
$ws.Range("J141").Value = 6298.6
$ws.Range("N141").Value = -29255.8
$ws.Range("L141").Value = 18895.8
$ws.Range("I141").Value = 2194.4211
$ws.Range("M141").Value = -1403.263300000001
$ws.Range("K141").Value = 6583.263300000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I2").Value = 1343.8
$ws.Range("K2").Value = 1343.8
$ws.Range("M2").Value = -1230.8
$ws.Range("H2").Value = 1314.7354
$ws.Range("M32").Value = -6456.277
$ws.Range("K32").Value = 6743.277
$ws.Range("H32").Value = 6723.2427
$ws.Range("I32").Value = 6743.277
$ws.Range("L32").Value = 6462.8
$ws.Range("N32").Value = -7036.8
$ws.Range("J32").Value = 6462.8
$ws.Range("M74").Value = -205.1538
$ws.Range("J74").Value = 6507
$ws.Range("L74").Value = 6507
$ws.Range("I74").Value = 1079.1538
$ws.Range("H74").Value = 3146.9048
$ws.Range("N74").Value = -8255
$ws.Range("K74").Value = 1079.1538
$ws.Range("M77").Value = -1027.769
$ws.Range("L77").Value = 32535
$ws.Range("N77").Value = -41271
$ws.Range("H77").Value = 3146.9048
$ws.Range("I77").Value = 1079.1538
$ws.Range("J77").Value = 6507
$ws.Range("K77").Value = 5395.769
$ws.Range("I116").Value = 1343.8
$ws.Range("M116").Value = 950.2
$ws.Range("K116").Value = 1343.8
$ws.Range("H116").Value = 1314.7354
$ws.Range("M132").Value = -4943.358200000001
$ws.Range("H132").Value = 2518.2764
$ws.Range("K132").Value = 7473.358200000001
$ws.Range("I132").Value = 2491.1194
$ws.Range("H139").Value = 91136.42999999999
$ws.Range("J139").Value = 91136.42999999999
$ws.Range("L139").Value = 91136.42999999999
$ws.Range("N139").Value = -101416.43

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M3").Value = -1229.8
$ws.Range("H3").Value = 1314.7354
$ws.Range("K3").Value = 1343.8
$ws.Range("I3").Value = 1343.8
$ws.Range("H5").Value = 1998.5
$ws.Range("I5").Value = 1998.5
$ws.Range("K5").Value = 1998.5
$ws.Range("M5").Value = -1885.5
$ws.Range("K99").Value = 2804.5
$ws.Range("H99").Value = 3376.6
$ws.Range("I99").Value = 2804.5
$ws.Range("M99").Value = -1306.5
$ws.Range("K107").Value = 1308.125
$ws.Range("H107").Value = 1661.5625
$ws.Range("J107").Value = 2015
$ws.Range("L107").Value = 2015
$ws.Range("M107").Value = 611.875
$ws.Range("N107").Value = -5855
$ws.Range("I107").Value = 1308.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I16").Value = 1675.5555
$ws.Range("K16").Value = 1675.5555
$ws.Range("M16").Value = -1388.5555
$ws.Range("H16").Value = 1908
$ws.Range("N31").Value = -5951.7144
$ws.Range("H31").Value = 2141.2
$ws.Range("K31").Value = 1336.0714
$ws.Range("M31").Value = -1041.0714
$ws.Range("I31").Value = 1336.0714
$ws.Range("J31").Value = 5361.7144
$ws.Range("L31").Value = 5361.7144
$ws.Range("H34").Value = 2141.2
$ws.Range("I34").Value = 1336.0714
$ws.Range("J34").Value = 5361.7144
$ws.Range("L34").Value = 5361.7144
$ws.Range("N34").Value = -5765.7144
$ws.Range("M34").Value = -1134.0714
$ws.Range("K34").Value = 1336.0714
$ws.Range("K113").Value = 1675.5555
$ws.Range("I113").Value = 1675.5555
$ws.Range("H113").Value = 1908
$ws.Range("M113").Value = 494.4445000000001
$ws.Range("M132").Value = -5003220.199999999
$ws.Range("H132").Value = 1252011.8
$ws.Range("K132").Value = 5005750.199999999
$ws.Range("I132").Value = 1668583.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M75").Value = -2762.0002
$ws.Range("I75").Value = 1253.3334
$ws.Range("H75").Value = 62109.176
$ws.Range("K75").Value = 3760.0002
$ws.Range("K78").Value = 11280.0006
$ws.Range("H78").Value = 62109.176
$ws.Range("I78").Value = 1253.3334
$ws.Range("M78").Value = -6288.000599999999
$ws.Range("M132").Value = -4808.2139
$ws.Range("H132").Value = 2267.8333
$ws.Range("K132").Value = 7338.2139
$ws.Range("I132").Value = 815.3570999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I44").Value = 0
$ws.Range("H44").Value = 25000
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("N80").Value = -6290.7
$ws.Range("J80").Value = 4294.7
$ws.Range("H80").Value = 4389.68
$ws.Range("L80").Value = 4294.7
$ws.Range("H83").Value = 4389.68
$ws.Range("J83").Value = 4294.7
$ws.Range("L83").Value = 21473.5
$ws.Range("N83").Value = -31457.5
$ws.Range("H122").Value = 3125.7036
$ws.Range("M122").Value = -3745.2001
$ws.Range("I122").Value = 2065.0667
$ws.Range("K122").Value = 6195.2001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("L20").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("L43").Value = 0
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H82").Value = 15202.637
$ws.Range("I82").Value = 27224.5
$ws.Range("K82").Value = 27224.5
$ws.Range("M82").Value = -26863.5
$ws.Range("I85").Value = 27224.5
$ws.Range("H85").Value = 15202.637
$ws.Range("M85").Value = -25976.5
$ws.Range("K85").Value = 27224.5
$ws.Range("M132").Value = -5926.163
$ws.Range("H132").Value = 2850.4106
$ws.Range("K132").Value = 8456.163
$ws.Range("I132").Value = 2818.721
$ws.Range("H134").Value = 78500
$ws.Range("H136").Value = 4760.625
$ws.Range("J136").Value = 5446.1
$ws.Range("L136").Value = 16338.3
$ws.Range("M136").Value = -10797.1362
$ws.Range("K136").Value = 13347.1362
$ws.Range("I136").Value = 4449.0454
$ws.Range("N136").Value = -21438.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K107").Value = 1274.7273
$ws.Range("H107").Value = 759.8421
$ws.Range("J107").Value = 1220.375
$ws.Range("L107").Value = 3661.125
$ws.Range("M107").Value = 645.2727
$ws.Range("N107").Value = -7501.125
$ws.Range("I107").Value = 424.9091
$ws.Range("N132").Value = -12286.5716
$ws.Range("J132").Value = 2408.8572
$ws.Range("M132").Value = -4681.000100000001
$ws.Range("H132").Value = 2405.12
$ws.Range("K132").Value = 7211.000100000001
$ws.Range("I132").Value = 2403.6667
$ws.Range("L132").Value = 7226.571599999999
$ws.Range("N133").Value = -122258
$ws.Range("L133").Value = 112138
$ws.Range("H133").Value = 112138
$ws.Range("J133").Value = 112138
